$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 5179.3335
$ws.Range("I74").Value = 4442.3335
$ws.Range("J74").Value = 5916.3335
$ws.Range("K74").Value = 4442.3335
$ws.Range("L74").Value = 5916.3335
$ws.Range("M74").Value = -3506.3335
$ws.Range("N74").Value = -7788.3335
# Row 77
$ws.Range("H77").Value = 5179.3335
$ws.Range("I77").Value = 4442.3335
$ws.Range("J77").Value = 5916.3335
$ws.Range("K77").Value = 22211.6675
$ws.Range("L77").Value = 29581.6675
$ws.Range("M77").Value = -17531.6675
$ws.Range("N77").Value = -38941.6675
# Row 80
$ws.Range("H80").Value = 6453.1055
$ws.Range("I80").Value = 666.6667
$ws.Range("J80").Value = 11660.9
$ws.Range("K80").Value = 2000.0001
$ws.Range("L80").Value = 34982.7
$ws.Range("M80").Value = -1002.0001
$ws.Range("N80").Value = -36978.7
# Row 83
$ws.Range("H83").Value = 6453.1055
$ws.Range("I83").Value = 666.6667
$ws.Range("J83").Value = 11660.9
$ws.Range("K83").Value = 6000.0003
$ws.Range("L83").Value = 104948.1
$ws.Range("M83").Value = -1008.0003
$ws.Range("N83").Value = -114932.1
# Row 135
$ws.Range("H135").Value = 62501584
$ws.Range("I135").Value = 23810810
$ws.Range("K135").Value = 214297290
$ws.Range("M135").Value = -214294755
# Row 138
$ws.Range("H138").Value = 15057.615
$ws.Range("I138").Value = 5259.2
$ws.Range("J138").Value = 17390.572
$ws.Range("K138").Value = 15777.6
$ws.Range("L138").Value = 52171.716
$ws.Range("M138").Value = -10637.6
$ws.Range("N138").Value = -62451.716

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 12177.485
$ws.Range("I61").Value = 9785.046
$ws.Range("J61").Value = 16226.23
$ws.Range("K61").Value = 9785.046
$ws.Range("L61").Value = 16226.23
$ws.Range("M61").Value = -9573.046
$ws.Range("N61").Value = -16650.23
# Row 63
$ws.Range("H63").Value = 3998.2
$ws.Range("I63").Value = 3492.5
$ws.Range("J63").Value = 4335.3335
$ws.Range("K63").Value = 3492.5
$ws.Range("L63").Value = 4335.3335
$ws.Range("M63").Value = -2806.5
$ws.Range("N63").Value = -5707.3335
# Row 66
$ws.Range("H66").Value = 3998.2
$ws.Range("I66").Value = 3492.5
$ws.Range("J66").Value = 4335.3335
$ws.Range("K66").Value = 17462.5
$ws.Range("L66").Value = 21676.6675
$ws.Range("M66").Value = -14030.5
$ws.Range("N66").Value = -28540.6675
# Row 136
$ws.Range("H136").Value = 12177.485
$ws.Range("I136").Value = 9785.046
$ws.Range("J136").Value = 16226.23
$ws.Range("K136").Value = 29355.138
$ws.Range("L136").Value = 48678.69
$ws.Range("M136").Value = -26805.138
$ws.Range("N136").Value = -53778.69

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1894.2593
$ws.Range("I86").Value = 1774.2
$ws.Range("J86").Value = 3395
$ws.Range("K86").Value = 1774.2
$ws.Range("L86").Value = 3395
$ws.Range("M86").Value = -651.2
$ws.Range("N86").Value = -5641
# Row 89
$ws.Range("H89").Value = 1894.2593
$ws.Range("I89").Value = 1774.2
$ws.Range("J89").Value = 3395
$ws.Range("K89").Value = 8871
$ws.Range("L89").Value = 16975
$ws.Range("M89").Value = -3255
$ws.Range("N89").Value = -28207

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 74
$ws.Range("H74").Value = 38875.668
$ws.Range("J74").Value = 38875.668
$ws.Range("L74").Value = 38875.668
$ws.Range("N74").Value = -40623.668
# Row 77
$ws.Range("H77").Value = 38875.668
$ws.Range("J77").Value = 38875.668
$ws.Range("L77").Value = 116627.004
$ws.Range("N77").Value = -125363.004
# Row 104
$ws.Range("H104").Value = 65000
$ws.Range("J104").Value = 65000
$ws.Range("L104").Value = 65000
$ws.Range("N104").Value = -70242

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 104
$ws.Range("H104").Value = 2042.125
$ws.Range("J104").Value = 2576.1667
$ws.Range("L104").Value = 7728.500100000001
$ws.Range("N104").Value = -12970.5001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 8223.333000000001
$ws.Range("I80").Value = 19066.666
$ws.Range("J80").Value = 4608.8887
$ws.Range("K80").Value = 19066.666
$ws.Range("L80").Value = 4608.8887
$ws.Range("M80").Value = -18068.666
$ws.Range("N80").Value = -6604.8887
# Row 83
$ws.Range("H83").Value = 8223.333000000001
$ws.Range("I83").Value = 19066.666
$ws.Range("J83").Value = 4608.8887
$ws.Range("K83").Value = 95333.33
$ws.Range("L83").Value = 23044.4435
$ws.Range("M83").Value = -90341.33
$ws.Range("N83").Value = -33028.4435
# Row 134
$ws.Range("H134").Value = 36550.668
$ws.Range("J134").Value = 36550.668
$ws.Range("L134").Value = 109652.004
$ws.Range("N134").Value = -114722.004

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3070.9443
$ws.Range("I40").Value = 2747.4
$ws.Range("J40").Value = 4688.6665
$ws.Range("K40").Value = 2747.4
$ws.Range("L40").Value = 4688.6665
$ws.Range("M40").Value = -2611.4
$ws.Range("N40").Value = -4960.6665
# Row 68
$ws.Range("H68").Value = 3110.4443
$ws.Range("I68").Value = 3124.25
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 3124.25
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -2375.25
$ws.Range("N68").Value = -4498
# Row 71
$ws.Range("H71").Value = 3110.4443
$ws.Range("I71").Value = 3124.25
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 15621.25
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -11877.25
$ws.Range("N71").Value = -22488
# Row 80
$ws.Range("H80").Value = 18000
$ws.Range("J80").Value = 18000
$ws.Range("L80").Value = 18000
$ws.Range("N80").Value = -20246
# Row 82
$ws.Range("H82").Value = 1966.5
$ws.Range("I82").Value = 1699.75
$ws.Range("J82").Value = 2500
$ws.Range("K82").Value = 1699.75
$ws.Range("L82").Value = 2500
$ws.Range("M82").Value = -1338.75
$ws.Range("N82").Value = -3222
# Row 83
$ws.Range("H83").Value = 18000
$ws.Range("J83").Value = 18000
$ws.Range("L83").Value = 54000
$ws.Range("N83").Value = -65232
# Row 85
$ws.Range("H85").Value = 1966.5
$ws.Range("I85").Value = 1699.75
$ws.Range("J85").Value = 2500
$ws.Range("K85").Value = 1699.75
$ws.Range("L85").Value = 2500
$ws.Range("M85").Value = -451.75
$ws.Range("N85").Value = -4996

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 2611
$ws.Range("I107").Value = 629
$ws.Range("J107").Value = 4052.4546
$ws.Range("K107").Value = 1887
$ws.Range("L107").Value = 12157.3638
$ws.Range("M107").Value = 33
$ws.Range("N107").Value = -15997.3638
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
# Row 111
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
# Row 112
$ws.Range("H112").Value = 72333.336
$ws.Range("J112").Value = 72333.336
$ws.Range("L112").Value = 72333.336
$ws.Range("N112").Value = -75287.336
# Row 113
$ws.Range("H113").Value = 743.1111
$ws.Range("I113").Value = 316.27274
$ws.Range("J113").Value = 1036.5625
$ws.Range("K113").Value = 948.81822
$ws.Range("L113").Value = 3109.6875
$ws.Range("M113").Value = 1221.18178
$ws.Range("N113").Value = -7449.6875
# Row 114
$ws.Range("H114").Value = 34800
$ws.Range("J114").Value = 34800
$ws.Range("L114").Value = 34800
$ws.Range("N114").Value = -43478
# Row 117
$ws.Range("H117").Value = 51566.668
$ws.Range("J117").Value = 51566.668
$ws.Range("L117").Value = 51566.668
$ws.Range("N117").Value = -60744.668

